# Apply the LOB1211 syllabus edit:
#  - Delete old row 13 (the stray "7455355 - Robson da Silva Rocha" value row with no A label)
#  - This shifts rows 14-22 up to become rows 13-21
#  - Update several B/C cell values to their new (shuffled) text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 13 entirely; rows below shift up by one.
$ws.Rows.Item(13).Delete()

# After the shift, fix up the cell contents that changed value.
$ws.Range("B10").Value = "7455355 - Robson da Silva Rocha"
$ws.Range("C10").Value = "7455355 - Robson da Silva Rocha"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B15").NumberFormat = "@"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("B15").Value = "01/01/2022"
$ws.Range("C15").Value = "01/01/2022"
$ws.Range("B15").NumberFormat = "General"
$ws.Range("C15").NumberFormat = "General"

$ws.Range("B18").Value = "7455355 - Robson da Silva Rocha"
$ws.Range("C18").Value = "7455355 - Robson da Silva Rocha"

$ws.Range("B19").Value = "Aulas teóricas expositivas e atividades em grupo."
$ws.Range("C19").Value = "Aulas teóricas expositivas e atividades em grupo."

$ws.Range("B20").Value = "Média ponderada de provas e atividades."
$ws.Range("C20").Value = "Média ponderada de provas e atividades."

$ws.Range("B21").Value = "1 (uma) prova escrita."
$ws.Range("C21").Value = "1 (uma) prova escrita."
